$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '''28.412.02'
$ws.Cells.Item(2, 5).Value = '''  +5.22%  '
$ws.Cells.Item(3, 4).Value = '''1.818.62'
$ws.Cells.Item(4, 4).Value = '''0.9988'
$ws.Cells.Item(4, 5).Value = '''  -0.25%  '
$ws.Cells.Item(5, 4).Value = '''317.98'
$ws.Cells.Item(5, 5).Value = '''  +2.39%  '
$ws.Cells.Item(6, 4).Value = '''0.9991'
$ws.Cells.Item(6, 5).Value = '''  -0.18%  '
$ws.Cells.Item(7, 4).Value = '''0.5773'
$ws.Cells.Item(7, 5).Value = '''  +18.99%  '
$ws.Cells.Item(8, 4).Value = '''0.3856'
$ws.Cells.Item(8, 5).Value = '''  +11.14%  '
$ws.Cells.Item(9, 4).Value = '''43.38'
$ws.Cells.Item(9, 5).Value = '''  +0.30%  '
$ws.Cells.Item(10, 4).Value = '''0.07636'
$ws.Cells.Item(10, 5).Value = '''  +5.62%  '
$ws.Cells.Item(11, 5).Value = '''  +8.62%  '
$ws.Cells.Item(12, 4).Value = '''21.32'
$ws.Cells.Item(12, 5).Value = '''  +7.05%  '
$ws.Cells.Item(13, 4).Value = '''0.9998'
$ws.Cells.Item(13, 5).Value = '''  -0.13%  '
$ws.Cells.Item(14, 5).Value = '''  +6.72%  '
$ws.Cells.Item(15, 4).Value = '''1.815.69'
$ws.Cells.Item(15, 5).Value = '''  +5.20%  '
$ws.Cells.Item(16, 4).Value = '''7.305'
$ws.Cells.Item(16, 5).Value = '''  +7.40%  '
$ws.Cells.Item(17, 4).Value = '''92.36'
$ws.Cells.Item(17, 5).Value = '''  +6.08%  '
$ws.Cells.Item(18, 4).Value = '''0.00001082'
$ws.Cells.Item(18, 5).Value = '''  +4.83%  '
$ws.Cells.Item(19, 4).Value = '''0.06523'
$ws.Cells.Item(19, 5).Value = '''  +1.92%  '
$ws.Cells.Item(20, 4).Value = '''0.9982'
$ws.Cells.Item(20, 5).Value = '''  -0.27%  '
$ws.Cells.Item(21, 4).Value = '''17.33'
$ws.Cells.Item(21, 5).Value = '''  +4.68%  '
$ws.Cells.Item(22, 4).Value = '''5.999'
$ws.Cells.Item(22, 5).Value = '''  +5.17%  '
$ws.Cells.Item(23, 4).Value = '''28.425.91'
$ws.Cells.Item(23, 5).Value = '''  +5.00%  '
$ws.Cells.Item(24, 4).Value = '''11.41'
$ws.Cells.Item(24, 5).Value = '''  +4.09%  '
$ws.Cells.Item(25, 4).Value = '''2.099'
$ws.Cells.Item(25, 5).Value = '''  +1.22%  '
$ws.Cells.Item(26, 4).Value = '''20.93'
$ws.Cells.Item(26, 5).Value = '''  +5.17%  '
$ws.Cells.Item(27, 4).Value = '''157.57'
$ws.Cells.Item(27, 5).Value = '''  +2.38%  '
$ws.Cells.Item(28, 4).Value = '''2.409'
$ws.Cells.Item(28, 5).Value = '''  +16.65%  '
$ws.Cells.Item(29, 4).Value = '''2.020.32'
$ws.Cells.Item(29, 5).Value = '''  +5.03%  '
$ws.Cells.Item(30, 4).Value = '''123.71'
$ws.Cells.Item(30, 5).Value = '''  +2.68%  '
$ws.Cells.Item(31, 4).Value = '''1.161'
$ws.Cells.Item(31, 5).Value = '''  +12.17%  '
$ws.Cells.Item(32, 4).Value = '''0.1061'
$ws.Cells.Item(32, 5).Value = '''  +13.86%  '
$ws.Cells.Item(33, 4).Value = '''5.780'
$ws.Cells.Item(33, 5).Value = '''  +7.65%  '
$ws.Cells.Item(34, 4).Value = '''3.625'
$ws.Cells.Item(34, 5).Value = '''  -0.28%  '
$ws.Cells.Item(35, 4).Value = '''0.02318'
$ws.Cells.Item(35, 5).Value = '''  +6.52%  '
$ws.Cells.Item(36, 2).Value = 'FraxShare'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(36, 4).Value = '''8.849'
$ws.Cells.Item(36, 5).Value = '''  +18.33%  '
$ws.Cells.Item(37, 2).Value = 'Algorand'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(37, 4).Value = '''0.2156'
$ws.Cells.Item(37, 5).Value = '''  +8.17%  '
$ws.Cells.Item(38, 4).Value = '''11.72'
$ws.Cells.Item(38, 5).Value = '''  +7.36%  '
$ws.Cells.Item(39, 4).Value = '''0.6439'
$ws.Cells.Item(39, 5).Value = '''  +7.92%  '
$ws.Cells.Item(40, 4).Value = '''0.06092'
$ws.Cells.Item(40, 5).Value = '''  +3.24%  '
$ws.Cells.Item(41, 4).Value = '''5.047'
$ws.Cells.Item(41, 5).Value = '''  +6.58%  '
$ws.Cells.Item(42, 2).Value = 'Frax'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(42, 4).Value = '''0.9985'
$ws.Cells.Item(42, 5).Value = '''  -0.19%  '
$ws.Cells.Item(43, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(43, 4).Value = '''1.159'
$ws.Cells.Item(43, 5).Value = '''  +3.72%  '
$ws.Cells.Item(44, 4).Value = '''1.380'
$ws.Cells.Item(44, 5).Value = '''  -3.55%  '
$ws.Cells.Item(45, 4).Value = '''13.48'
$ws.Cells.Item(46, 4).Value = '''0.5993'
$ws.Cells.Item(46, 5).Value = '''  +7.01%  '
$ws.Cells.Item(47, 4).Value = '''3.706'
$ws.Cells.Item(47, 5).Value = '''  +3.48%  '
$ws.Cells.Item(48, 4).Value = '''122.25'
$ws.Cells.Item(48, 5).Value = '''  +2.52%  '
$ws.Cells.Item(49, 4).Value = '''1.946'
$ws.Cells.Item(49, 5).Value = '''  +5.65%  '
$ws.Cells.Item(50, 4).Value = '''1.146'
$ws.Cells.Item(50, 5).Value = '''  +4.12%  '
$ws.Cells.Item(51, 4).Value = '''0.06854'
